# This workbook contains a weekly price log. The update inserts one new
# weekly record at row 475 (pushing the existing rows 475-566 down to
# 476-567) and populates the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 475, shifting existing data down.
$ws.Rows("475:475").Insert()

# Populate the new row 475 with the new weekly record.
$ws.Cells.Item(475, 1).Value = 4
$ws.Cells.Item(475, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(475, 3).Value = "Los Lagos"
$ws.Cells.Item(475, 4).Value = 45258
$ws.Cells.Item(475, 5).Value = 10
$ws.Cells.Item(475, 6).Value = 100114014
$ws.Cells.Item(475, 7).Value = "Betarraga"
$ws.Cells.Item(475, 8).Value = "Sin especificar"
$ws.Cells.Item(475, 9).Value = "Primera"
$ws.Cells.Item(475, 10).Value = 1200
$ws.Cells.Item(475, 11).Value = 1000
$ws.Cells.Item(475, 12).Value = 1100
$ws.Cells.Item(475, 13).Value = 1050
$ws.Cells.Item(475, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(475, 15).Value = "Región Metropolitana"
$ws.Cells.Item(475, 16).Value = 210
$ws.Cells.Item(475, 17).Value = 5
$ws.Cells.Item(475, 18).Value = "Hortaliza"
